$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 1.9
$ws.Range("Q3").Value = 2
$ws.Range("G4").Value = 1.47
$ws.Range("H4").Value = 10
$ws.Range("H5").Value = 3.1
$ws.Range("I5").Value = 4.5
$ws.Range("P7").Value = 2.38
$ws.Range("Q7").Value = 1.47
$ws.Range("F8").Value = 4.5
$ws.Range("I8").Value = 1.85
$ws.Range("J8").Value = 4.2
$ws.Range("Q8").Value = 1.79
$ws.Range("T8").Value = 1.78
$ws.Range("U8").Value = 2.2
$ws.Range("X8").Value = 19
$ws.Range("AA8").Value = 19.5
$ws.Range("AB8").Value = 18.5
$ws.Range("AF8").Value = 38
$ws.Range("AG8").Value = 18.5
$ws.Range("AH8").Value = 19
$ws.Range("AJ8").Value = 130
$ws.Range("AK8").Value = 60
$ws.Range("AL8").Value = 60
$ws.Range("AN8").Value = 60
$ws.Range("AO8").Value = 10
$ws.Range("H10").Value = 2.32
$ws.Range("I10").Value = 2.34
$ws.Range("Q10").Value = 2.18
$ws.Range("AM10").Value = 130
$ws.Range("F11").Value = 1.34
$ws.Range("G11").Value = 1.35
$ws.Range("H11").Value = 11.5
$ws.Range("K11").Value = 5.9
$ws.Range("S11").Value = 2.92
$ws.Range("T11").Value = 2.24
$ws.Range("Z11").Value = 110
$ws.Range("AA11").Value = 510
$ws.Range("AF11").Value = 7.6
$ws.Range("AM11").Value = 250
$ws.Range("F12").Value = 4.8
$ws.Range("G12").Value = 4.9
$ws.Range("H12").Value = 1.84
$ws.Range("I12").Value = 1.86
$ws.Range("J12").Value = 3.9
$ws.Range("K12").Value = 3.95
$ws.Range("L12").Value = 1.08
$ws.Range("T12").Value = 1.96
$ws.Range("U12").Value = 1.97
$ws.Range("V12").Value = 2.16
$ws.Range("W12").Value = 1.25
$ws.Range("Y12").Value = 8.4
$ws.Range("Z12").Value = 11
$ws.Range("AA12").Value = 19.5
$ws.Range("AB12").Value = 15.5
$ws.Range("AE12").Value = 21
$ws.Range("AJ12").Value = 140
$ws.Range("AK12").Value = 75
$ws.Range("AN12").Value = 100
$ws.Range("AO12").Value = 14
$ws.Range("F13").Value = 2.22
$ws.Range("I13").Value = 3.6
$ws.Range("P13").Value = 2.06
$ws.Range("W13").Value = 1.8
$ws.Range("Y13").Value = 14.5
$ws.Range("G14").Value = 2.9
$ws.Range("Q14").Value = 1.76
$ws.Range("I15").Value = 1.78
$ws.Range("Q15").Value = 1.6
$ws.Range("T15").Value = 1.71
$ws.Range("U15").Value = 2.1
$ws.Range("AN15").Value = 65
$ws.Range("G16").Value = 1.39
$ws.Range("J16").Value = 5.3
$ws.Range("K16").Value = 5.9
$ws.Range("L16").Value = 1.28
$ws.Range("P16").Value = 2.22
$ws.Range("Q16").Value = 1.71
$ws.Range("R16").Value = 1.49
$ws.Range("S16").Value = 2.78
$ws.Range("U16").Value = 1.76
$ws.Range("W16").Value = 3.5
$ws.Range("X16").Value = 25
$ws.Range("AC16").Value = 13
$ws.Range("AD16").Value = 1000
$ws.Range("AF16").Value = 8.4
$ws.Range("AN16").Value = 6.8
$ws.Range("G17").Value = 1.12
$ws.Range("H17").Value = 27
$ws.Range("I17").Value = 36
$ws.Range("J17").Value = 12.5
$ws.Range("K17").Value = 15
$ws.Range("L17").Value = 1.14
$ws.Range("N17").Value = 8.6
$ws.Range("O17").Value = 1.09
$ws.Range("P17").Value = 3.95
$ws.Range("Q17").Value = 1.27
$ws.Range("R17").Value = 2.18
$ws.Range("S17").Value = 1.67
$ws.Range("T17").Value = 2.18
$ws.Range("U17").Value = 1.68
$ws.Range("V17").Value = 1.02
$ws.Range("W17").Value = 9.199999999999999
$ws.Range("AB17").Value = 17
$ws.Range("AC17").Value = 34
$ws.Range("AF17").Value = 11
$ws.Range("AG17").Value = 16.5
$ws.Range("AJ17").Value = 9.4
$ws.Range("AK17").Value = 18
$ws.Range("AN17").Value = 2.4
$ws.Range("G18").Value = 2.48
$ws.Range("H18").Value = 3.15
$ws.Range("I18").Value = 3.55
$ws.Range("J18").Value = 3.35
$ws.Range("L18").Value = 1.41
$ws.Range("M18").Value = 1.07
$ws.Range("N18").Value = 3.65
$ws.Range("O18").Value = 1.31
$ws.Range("P18").Value = 1.91
$ws.Range("Q18").Value = 1.91
$ws.Range("R18").Value = 1.35
$ws.Range("S18").Value = 3.3
$ws.Range("T18").Value = 1.72
$ws.Range("U18").Value = 2.16
$ws.Range("V18").Value = 1.39
$ws.Range("W18").Value = 1.67
$ws.Range("X18").Value = 17.5
$ws.Range("Y18").Value = 14
$ws.Range("AB18").Value = 11
$ws.Range("AC18").Value = 8.199999999999999
$ws.Range("AD18").Value = 15
$ws.Range("AF18").Value = 16
$ws.Range("AG18").Value = 11.5
$ws.Range("AH18").Value = 17.5
$ws.Range("I19").Value = 9.4
$ws.Range("Q19").Value = 1.62
$ws.Range("T19").Value = 1.94
$ws.Range("AI19").Value = 970
$ws.Range("AM19").Value = 120
$ws.Range("I20").Value = 1.84
$ws.Range("J20").Value = 4.1
$ws.Range("K20").Value = 4.2
$ws.Range("L20").Value = 1.31
$ws.Range("N20").Value = 4.1
$ws.Range("P20").Value = 2.08
$ws.Range("Q20").Value = 1.87
$ws.Range("S20").Value = 3.3
$ws.Range("T20").Value = 1.84
$ws.Range("U20").Value = 2.1
$ws.Range("V20").Value = 2.18
$ws.Range("Y20").Value = 9.4
$ws.Range("AD20").Value = 9.800000000000001
$ws.Range("AK20").Value = 60
$ws.Range("AL20").Value = 65
$ws.Range("AM20").Value = 110
$ws.Range("F21").Value = 2.02
$ws.Range("K21").Value = 3.9
$ws.Range("F22").Value = 2.34
$ws.Range("G22").Value = 2.54
$ws.Range("J22").Value = 3.4
$ws.Range("P22").Value = 1.81
$ws.Range("G23").Value = 5.7
$ws.Range("H23").Value = 1.75
$ws.Range("I23").Value = 1.87
$ws.Range("J23").Value = 3.6
$ws.Range("Q23").Value = 1.98
